# Finished simulation portion of report
#
# Fill in the "200k" sheet (4th sheet, previously blank) with the
# Analysis / Simulation / Experimentation comparison table, matching the
# style already used on the "150k" and "250k" sheets, and make it the
# active/selected sheet since it's the one most recently worked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("200k")

# --- Column widths (matches the layout used for the other comparison sheets) ---
$ws.Columns.Item(1).ColumnWidth = 9.140625
$ws.Columns.Item(2).ColumnWidth = 9.28515625
$ws.Columns.Item(3).ColumnWidth = 9.5703125
$ws.Columns.Item(4).ColumnWidth = 9.28515625
$ws.Columns.Item(5).ColumnWidth = 17.42578125
$ws.Columns.Item(6).ColumnWidth = 9.5703125
$ws.Columns.Item(7).ColumnWidth = 9.28515625

# --- Formatting: left-aligned, 5-decimal custom number format for the table ---
$ws.Range("B1:J2,A3:G9").NumberFormat = "0.00000"
$ws.Range("B1:J2,A3:G9").HorizontalAlignment = -4131

# --- Header row 1: merged group headers ---
$ws.Range("B1:D1").Merge()
$ws.Range("B1").Value = "Analysis"

$ws.Range("E1:G1").Merge()
$ws.Range("E1").Value = "Simulation"

$ws.Range("H1:J1").Merge()
$ws.Range("H1").Value = "Experimentation"

# --- Header row 2: column labels ---
$ws.Range("B2").Value = "t (ms)"
$ws.Range("C2").Value = "vin (V)"
$ws.Range("D2").Value = "vout (V)"
$ws.Range("E2").Value = "t (ms)"
$ws.Range("F2").Value = "vin (V)"
$ws.Range("G2").Value = "vout (V)"
$ws.Range("H2").Value = "t (ms)"
$ws.Range("I2").Value = "vin (V)"
$ws.Range("J2").Value = "vout (V)"

# --- Data rows 3-9: row labels ---
$ws.Range("A3").Value = "Breakover"
$ws.Range("A4").Value = "Peak 1"
$ws.Range("A5").Value = "Peak 2"
$ws.Range("A6").Value = "Peak 3"
$ws.Range("A7").Value = "Peak 4"
$ws.Range("A8").Value = "Peak 5"
$ws.Range("A9").Value = "Peak 6"

# --- Data rows 3-9: Analysis (B:D) / Simulation (E:G) values ---
$ws.Range("B3").Value = 0.00589
$ws.Range("C3").Value = 135.36702102541784
$ws.Range("D3").Value = 32.01848344222199
$ws.Range("E3").Value = 0.005931960810184
$ws.Range("F3").Value = 133.7233643128
$ws.Range("G3").Value = 32.10880887227

$ws.Range("B4").Value = 0.01621
$ws.Range("C4").Value = 28.912224463505513
$ws.Range("D4").Value = 28.54618047946137
$ws.Range("E4").Value = 0.01625381170887
$ws.Range("F4").Value = 26.3525525907
$ws.Range("G4").Value = 28.39757540369

$ws.Range("B5").Value = 0.03294
$ws.Range("C5").Value = 24.90455819719604
$ws.Range("D5").Value = 24.958450651225178
$ws.Range("E5").Value = 0.03297262527834
$ws.Range("F5").Value = 23.0460547843
$ws.Range("G5").Value = 24.82900158544

$ws.Range("B6").Value = 0.04963
$ws.Range("C6").Value = 23.42434941638862
$ws.Range("D6").Value = 23.391896131682675
$ws.Range("E6").Value = 0.04966803538733
$ws.Range("F6").Value = 21.2196193251
$ws.Range("G6").Value = 23.2829856851

$ws.Range("B7").Value = 0.06631
$ws.Range("C7").Value = 22.577694694958705
$ws.Range("D7").Value = 22.70709665741043
$ws.Range("E7").Value = 0.06628180997305
$ws.Range("F7").Value = 24.57844116957
$ws.Range("G7").Value = 22.61355985712

$ws.Range("B8").Value = 0.08298
$ws.Range("C8").Value = 22.365941045689226
$ws.Range("D8").Value = 22.407608274481422
$ws.Range("E8").Value = 0.08299199037025
$ws.Range("F8").Value = 21.81582915381
$ws.Range("G8").Value = 22.32547890283

$ws.Range("B9").Value = 0.09965
$ws.Range("C9").Value = 22.154152077545362
$ws.Range("D9").Value = 22.276591015932237
$ws.Range("E9").Value = 0.0996753097352
$ws.Range("F9").Value = 20.75698419488
$ws.Range("G9").Value = 22.1987763395

# --- Selection state on the previously-active sheet (250k) changed too ---
$ws3 = $wb.Worksheets.Item("250k")
$ws3.Range("A3:A8").Select()

# --- Make the new sheet (200k) the active tab, with B10 selected ---
$ws.Activate()
$ws.Range("B10").Select()
